$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1): F2 314 -> 315, F4 49 -> 51
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 315
$ws1.Range("F4").Value = 51

# Update "全部类型" sheet (sheet4): F2 314 -> 315, F4 49 -> 51
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 315
$ws4.Range("F4").Value = 51
